$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 10 data: B10 label + C10 numeric speed value
$ws.Range("B10").Value = "Why with decorator the time is small (in all cases)"
$ws.Range("C10").Value = 0.001

# Widen column B to fit the new text
$ws.Columns.Item(2).ColumnWidth = 44

# Move the active cell selection to C12
$ws.Range("C12").Select()
